# "Final Detection fully validated"
# Updates the detection-validation sheet: refreshes a handful of corrected
# STOP/CEDA detection counts, clears the not-yet-validated DETECTED_S /
# DETECTEC_C counts for the remaining rows, and nudges the column widths
# (A-E) to their final autofit sizes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column width adjustments (closest achievable values) ---
$ws.Columns.Item(1).ColumnWidth = 39.75
$ws.Columns.Item(2).ColumnWidth = 4.25
$ws.Columns.Item(3).ColumnWidth = 4.583333333333333
$ws.Columns.Item(4).ColumnWidth = 10.583333333333334
$ws.Columns.Item(5).ColumnWidth = 10.583333333333334

# --- Specific cell value corrections (DETECTED_S / DETECTEC_C columns) ---
$ws.Range("D3").Value = 1
$ws.Range("D8").Value = 0
$ws.Range("E10").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 1
$ws.Range("E18").Value = 1
$ws.Range("D19").Value = 1
$ws.Range("E20").Value = 2

# --- Clear DETECTED_S / DETECTEC_C columns for rows 22-95 (not yet validated) ---
$clearRange = $ws.Range("D22:E95")
$clearRange.ClearContents()
$clearRange.Style = "Normal"
